$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 16:18. This removes 3 rows and shifts the old row 19
# (Reebok CL Lthr Mid GTX-Thin) up into row 16's place, while the old
# rows 16-18 are discarded. The net result matches the target diff:
# old rows 16,17,18 are gone and old row 19's content now lives at row 16.
$ws.Range("A16:D18").EntireRow.Delete()
